# Scheduled runner update: refresh computed profit-sheet values (currentAveragePrice*,
# LevePrice*, LeveProfit* columns H-N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 591.8333
$ws.Range("J41").Value = 533.3333
$ws.Range("L41").Value = 533.3333
$ws.Range("N41").Value = -1413.3333

$ws.Range("H76").Value = 3520.8
$ws.Range("J76").Value = 3769.3333
$ws.Range("L76").Value = 3769.3333
$ws.Range("N76").Value = -4399.3333

$ws.Range("H79").Value = 3520.8
$ws.Range("J79").Value = 3769.3333
$ws.Range("L79").Value = 3769.3333
$ws.Range("N79").Value = -5953.3333

$ws.Range("H125").Value = 1722.1818
$ws.Range("I125").Value = 916
$ws.Range("J125").Value = 1901.3334
$ws.Range("K125").Value = 8244
$ws.Range("L125").Value = 17112.0006
$ws.Range("M125").Value = -5784
$ws.Range("N125").Value = -22032.0006

$ws.Range("H132").Value = 1169100.1
$ws.Range("I132").Value = 1914.2307
$ws.Range("J132").Value = 3065777.2
$ws.Range("K132").Value = 5742.6921
$ws.Range("L132").Value = 9197331.600000001
$ws.Range("M132").Value = -3212.6921
$ws.Range("N132").Value = -9202391.600000001

$ws.Range("H137").Value = 7696364
$ws.Range("I137").Value = 25002944
$ws.Range("J137").Value = 4550.3335
$ws.Range("K137").Value = 75008832
$ws.Range("L137").Value = 13651.0005
$ws.Range("M137").Value = -75006282
$ws.Range("N137").Value = -18751.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1295.5588
$ws.Range("I2").Value = 1193.2858
$ws.Range("J2").Value = 1460.7693
$ws.Range("K2").Value = 1193.2858
$ws.Range("L2").Value = 1460.7693
$ws.Range("M2").Value = -1080.2858
$ws.Range("N2").Value = -1686.7693

$ws.Range("H32").Value = 6572.77
$ws.Range("I32").Value = 3449.074
$ws.Range("J32").Value = 19889.578
$ws.Range("K32").Value = 3449.074
$ws.Range("L32").Value = 19889.578
$ws.Range("M32").Value = -3162.074
$ws.Range("N32").Value = -20463.578

$ws.Range("H45").Value = 1400
$ws.Range("I45").Value = 1600
$ws.Range("J45").Value = 1100
$ws.Range("K45").Value = 1600
$ws.Range("L45").Value = 1100
$ws.Range("M45").Value = -1223
$ws.Range("N45").Value = -1854

$ws.Range("H61").Value = 33402220
$ws.Range("I61").Value = 50051480
$ws.Range("J61").Value = 103698.7
$ws.Range("K61").Value = 50051480
$ws.Range("L61").Value = 103698.7
$ws.Range("M61").Value = -50051268
$ws.Range("N61").Value = -104122.7

$ws.Range("H74").Value = 8180273.5
$ws.Range("I74").Value = 10132243
$ws.Range("J74").Value = 128401.25
$ws.Range("K74").Value = 10132243
$ws.Range("L74").Value = 128401.25
$ws.Range("M74").Value = -10131369
$ws.Range("N74").Value = -130149.25

$ws.Range("H77").Value = 8180273.5
$ws.Range("I77").Value = 10132243
$ws.Range("J77").Value = 128401.25
$ws.Range("K77").Value = 50661215
$ws.Range("L77").Value = 642006.25
$ws.Range("M77").Value = -50656847
$ws.Range("N77").Value = -650742.25

$ws.Range("H116").Value = 1295.5588
$ws.Range("I116").Value = 1193.2858
$ws.Range("J116").Value = 1460.7693
$ws.Range("K116").Value = 1193.2858
$ws.Range("L116").Value = 1460.7693
$ws.Range("M116").Value = 1100.7142
$ws.Range("N116").Value = -6048.7693

$ws.Range("H122").Value = 6537720
$ws.Range("I122").Value = 1351.6154
$ws.Range("K122").Value = 4054.8462
$ws.Range("M122").Value = -1604.8462

$ws.Range("H132").Value = 72337.766
$ws.Range("I132").Value = 44193.207
$ws.Range("J132").Value = 184916
$ws.Range("K132").Value = 132579.621
$ws.Range("L132").Value = 554748
$ws.Range("M132").Value = -130049.621
$ws.Range("N132").Value = -559808

$ws.Range("H136").Value = 33402220
$ws.Range("I136").Value = 50051480
$ws.Range("J136").Value = 103698.7
$ws.Range("K136").Value = 150154440
$ws.Range("L136").Value = 311096.1
$ws.Range("M136").Value = -150151890
$ws.Range("N136").Value = -316196.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1295.5588
$ws.Range("I3").Value = 1193.2858
$ws.Range("J3").Value = 1460.7693
$ws.Range("K3").Value = 1193.2858
$ws.Range("L3").Value = 1460.7693
$ws.Range("M3").Value = -1079.2858
$ws.Range("N3").Value = -1688.7693

$ws.Range("H105").Value = 41668790
$ws.Range("I105").Value = 55557656
$ws.Range("J105").Value = 2200
$ws.Range("K105").Value = 55557656
$ws.Range("L105").Value = 2200
$ws.Range("M105").Value = -55555909
$ws.Range("N105").Value = -5694

$ws.Range("H134").Value = 2717.7334
$ws.Range("I134").Value = 2717.7334
$ws.Range("K134").Value = 8153.2002
$ws.Range("M134").Value = -5618.2002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 157.5625
$ws.Range("I7").Value = 74
$ws.Range("J7").Value = 296.83334
$ws.Range("K7").Value = 74
$ws.Range("L7").Value = 296.83334
$ws.Range("M7").Value = 39
$ws.Range("N7").Value = -522.83334

$ws.Range("H31").Value = 2802.3428
$ws.Range("I31").Value = 2384.5
$ws.Range("K31").Value = 2384.5
$ws.Range("M31").Value = -2089.5

$ws.Range("H34").Value = 2802.3428
$ws.Range("I34").Value = 2384.5
$ws.Range("K34").Value = 2384.5
$ws.Range("M34").Value = -2182.5

$ws.Range("H58").Value = 45456720
$ws.Range("I58").Value = 66668596
$ws.Range("J58").Value = 2703.4285
$ws.Range("K58").Value = 66668596
$ws.Range("L58").Value = 2703.4285
$ws.Range("M58").Value = -66668393
$ws.Range("N58").Value = -3109.4285

$ws.Range("H99").Value = 1151.7142
$ws.Range("I99").Value = 1120.3636
$ws.Range("J99").Value = 1266.6666
$ws.Range("K99").Value = 1120.3636
$ws.Range("L99").Value = 1266.6666
$ws.Range("M99").Value = 377.6364000000001
$ws.Range("N99").Value = -4262.6666

$ws.Range("H126").Value = 1151.7142
$ws.Range("I126").Value = 1120.3636
$ws.Range("J126").Value = 1266.6666
$ws.Range("K126").Value = 3361.0908
$ws.Range("L126").Value = 3799.9998
$ws.Range("M126").Value = -891.0907999999999
$ws.Range("N126").Value = -8739.9998

$ws.Range("H132").Value = 20369.537
$ws.Range("I132").Value = 1479.2188
$ws.Range("J132").Value = 47846.363
$ws.Range("K132").Value = 4437.6564
$ws.Range("L132").Value = 143539.089
$ws.Range("M132").Value = -1907.6564
$ws.Range("N132").Value = -148599.089

$ws.Range("H134").Value = 61972.89
$ws.Range("I134").Value = 1802.4
$ws.Range("J134").Value = 85115.38
$ws.Range("K134").Value = 5407.200000000001
$ws.Range("L134").Value = 255346.14
$ws.Range("M134").Value = -2872.200000000001
$ws.Range("N134").Value = -260416.14

$ws.Range("H136").Value = 45456720
$ws.Range("I136").Value = 66668596
$ws.Range("J136").Value = 2703.4285
$ws.Range("K136").Value = 200005788
$ws.Range("L136").Value = 8110.2855
$ws.Range("M136").Value = -200003238
$ws.Range("N136").Value = -13210.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 675.55554
$ws.Range("I113").Value = 590
$ws.Range("J113").Value = 680.58826
$ws.Range("K113").Value = 1770
$ws.Range("L113").Value = 2041.76478
$ws.Range("M113").Value = 400
$ws.Range("N113").Value = -6381.76478

$ws.Range("H134").Value = 4394.6
$ws.Range("I134").Value = 2479.6875
$ws.Range("J134").Value = 7798.8887
$ws.Range("K134").Value = 7439.0625
$ws.Range("L134").Value = 23396.6661
$ws.Range("M134").Value = -2369.0625
$ws.Range("N134").Value = -33536.6661

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 69.666664
$ws.Range("I2").Value = 9.25
$ws.Range("J2").Value = 91.63636
$ws.Range("K2").Value = 9.25
$ws.Range("L2").Value = 91.63636
$ws.Range("M2").Value = 103.75
$ws.Range("N2").Value = -317.63636

$ws.Range("H126").Value = 1819.7
$ws.Range("I126").Value = 1504
$ws.Range("J126").Value = 1955
$ws.Range("K126").Value = 4512
$ws.Range("L126").Value = 5865
$ws.Range("M126").Value = -2042
$ws.Range("N126").Value = -10805

$ws.Range("H132").Value = 336240.34
$ws.Range("I132").Value = 1000000
$ws.Range("J132").Value = 203488.4
$ws.Range("K132").Value = 3000000
$ws.Range("L132").Value = 610465.2
$ws.Range("M132").Value = -2997470
$ws.Range("N132").Value = -615525.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 98986.09
$ws.Range("I136").Value = 62335.47
$ws.Range("J136").Value = 254751.25
$ws.Range("K136").Value = 187006.41
$ws.Range("L136").Value = 764253.75
$ws.Range("M136").Value = -184456.41
$ws.Range("N136").Value = -769353.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1010.4167
$ws.Range("I113").Value = 512.5
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 1537.5
$ws.Range("L113").Value = 10500
$ws.Range("M113").Value = 632.5
$ws.Range("N113").Value = -14840

$ws.Range("H132").Value = 58896.113
$ws.Range("I132").Value = 53981.633
$ws.Range("K132").Value = 161944.899
$ws.Range("M132").Value = -159414.899
